$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of old single-letter CAENES codes to new merged-section codes
$map = @{
    "S" = "RS"
    "D" = "DE"
    "H" = "HJ"
    "L" = "LMN"
    "R" = "RS"
}

$lastRow = 171
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}

# Update the view state to match the saved workbook (scrolled/selected state)
$ws.Application.ActiveWindow.ScrollRow = 156
$ws.Range("D1:D1048576").Select()
